# Apply TabuSearch_Stats.xlsx data updates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 - U7
$ws.Range("B2").Value = 526.45
$ws.Range("C2").Value = 195.0
$ws.Range("E2").Value = 185.0

# Row 3 - U8
$ws.Range("B3").Value = 725.27
$ws.Range("C3").Value = 204.0
$ws.Range("E3").Value = 202.0

# Row 4 - U9
$ws.Range("B4").Value = 834.11
$ws.Range("C4").Value = 344.0
$ws.Range("D4").Value = 14.0
$ws.Range("E4").Value = 330.0

# Row 16 - Division U7 (tier: 3)
$ws.Range("B16").Value = 324.75
$ws.Range("C16").Value = 1723.0
$ws.Range("D16").Value = 35.0
$ws.Range("E16").Value = 1688.0

# Row 17 - Division U8 (tier: 3)
$ws.Range("B17").Value = 425.27
$ws.Range("C17").Value = 990.0
$ws.Range("D17").Value = 50.0
$ws.Range("E17").Value = 940.0

# Row 18 - Division U9 (tier: 3)
$ws.Range("B18").Value = 736.11
$ws.Range("C18").Value = 1820.0
$ws.Range("D18").Value = 60.0
$ws.Range("E18").Value = 1760.0

# Row 30 - Entire League
$ws.Range("B30").Value = 1304.87
$ws.Range("C30").Value = 9112.0
$ws.Range("D30").Value = 233.0
$ws.Range("E30").Value = 8879.0
$ws.Range("F30").Value = "0 min, 10 sec"
